# Apply the "added graph to population page; resize works" edit.
#
# Summary of the change (derived from the OOXML diff):
#  - A new (blank) column is inserted at column C on Sheet1, pushing the
#    old column C -> D and old column D -> E.
#  - The header row gets a new label in the freshly inserted C1: "var name".
#  - A brand-new trailing data row (row 97) is appended describing a new
#    "Analysis" / "interacts with the freezer" entry in columns A and D.
#  - The sheet view scroll/selection resets to the top of the sheet with
#    C2 selected (previously it was scrolled down to row 61 with C74
#    selected).
#  - The workbook window's on-screen position changes (xWindow/yWindow).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column C. Excel's native "insert column"
# behavior shifts the old C/D columns right one slot and copies the
# formatting of the inserted cells from the existing row styles, which is
# exactly the pattern seen throughout the diff (each shifted row gains an
# empty, styled C cell before the relocated D cell).
$ws.Columns("C").Insert()

# Header row: label the newly inserted column.
$ws.Range("C1").Value = "var name"

# New trailing row describing the freezer-interaction analysis entry.
$ws.Range("A97").Value = "Analysis"
$ws.Range("D97").Value = "interacts with the freezer"

# Reset the view: scroll back to the top and select C2 (was topLeftCell
# A61 / selection C74 before the edit).
$ws.Range("C2").Select()

# Move the workbook window on screen to match the new recorded position.
$win = $excel.ActiveWindow
$win.Left = 28680
$win.Top = 60
